$d = $word.ActiveDocument

# 1. Plea: "No Contest" -> "Guilty"
$d.Content.Find.Execute("No Contest", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Guilty", 2)

# 2. Community control period: "1 year" -> "18 months"
$d.Content.Find.Execute("1 year", $true, $false, $false, $false, $false,
                         $true, 1, $false, "18 months", 2)

# 3. Remove the two "no contact"/"stay away" conditions paragraphs that were
#    refactored out of this document (condition terms moved elsewhere).
#    Locate them by their distinctive text.
$n = $d.Paragraphs.Count
$idxFeet = -1
$idxRed = -1
$idxEmployment = -1
for ($i = 1; $i -le $n; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*500 feet*") {
        $idxFeet = $i
    }
    elseif ($t -like "*Defendant shall have no contact with*") {
        $idxRed = $i
    }
    elseif ($t -like "*Make all reasonable effort to obtain and maintain employment*") {
        $idxEmployment = $i
    }
}

# Drop the trailing standalone space run left at the end of the
# "Make all reasonable effort..." paragraph (the run immediately preceding
# the two paragraphs being removed).
$pEmployment = $d.Paragraphs($idxEmployment)
$trimStart = $pEmployment.Range.End - 2
$trimEnd = $pEmployment.Range.End - 1
$d.Range($trimStart, $trimEnd).Delete()

# Delete both paragraphs (and their paragraph marks) in one shot.
$pFeet = $d.Paragraphs($idxFeet)
$pRed = $d.Paragraphs($idxRed)
$d.Range($pFeet.Range.Start, $pRed.Range.End).Delete()
